$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0
$ws.Range("B2").Value = -0.1056763790825679
$ws.Range("C2").Value = -0
$ws.Range("D2").Value = 0.2661582659571812
$ws.Range("E2").Value = 0.01003772731320099
$ws.Range("F2").Value = -0
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = -0
$ws.Range("J2").Value = -0
$ws.Range("K2").Value = -0.01979208801498059
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.242987368219968
$ws.Range("N2").Value = 0.004423945506101552
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.1076971022137111
$ws.Range("V2").Value = 0.01605192199234246
$ws.Range("W2").Value = -0.04460431739448262
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = -0.07772532331582976
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = -0.00988950490464704
$ws.Range("AF2").Value = 0.009194210549037202
$ws.Range("AG2").Value = -0
$ws.Range("AI2").Value = -0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = -0.04384480141443917
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.03442446281322958
$ws.Range("AO2").Value = 0.07982703158461509
$ws.Range("AQ2").Value = 0
$ws.Range("AR2").Value = -0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1648594835943316
$ws.Range("AW2").Value = 0.09513739437549891
$ws.Range("AX2").Value = 0.01400073286441289
$ws.Range("AY2").Value = -0
$ws.Range("BB2").Value = -0
$ws.Range("BC2").Value = -0
$ws.Range("BD2").Value = -0.02905379851713809
$ws.Range("BF2").Value = 0.09983054286285226
$ws.Range("BG2").Value = 0.03683329100565719
$ws.Range("BJ2").Value = -0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.02957401969734591
$ws.Range("BO2").Value = -0.04199075926153524
$ws.Range("BP2").Value = -0.09521321784165342
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.04027809927706854
$ws.Range("BX2").Value = 0.02140014096010549
$ws.Range("BY2").Value = -0.02033857839823998
$ws.Range("BZ2").Value = -0
$ws.Range("CB2").Value = 0
$ws.Range("CD2").Value = -0
$ws.Range("CE2").Value = 0.03341464267615306
$ws.Range("CG2").Value = -0.04024297188919205
$ws.Range("CH2").Value = 0.01772990530384111
$ws.Range("CJ2").Value = -0
$ws.Range("CM2").Value = -0
$ws.Range("CN2").Value = -0.01476736764181623
$ws.Range("CP2").Value = 0.03804480328825097
$ws.Range("CQ2").Value = 0.04306717921016637
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = -0
$ws.Range("CV2").Value = -0
$ws.Range("CW2").Value = 0.04779099857713974
$ws.Range("CY2").Value = -0.04395842935551074
$ws.Range("CZ2").Value = 0.009211424765332644
$ws.Range("DE2").Value = -0
$ws.Range("DF2").Value = 0.03461260656283643
$ws.Range("DH2").Value = 0.01670588914064872
$ws.Range("DI2").Value = 0.0364763747653508
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = -0
$ws.Range("DL2").Value = -0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.03330360949037436
$ws.Range("DQ2").Value = 0.04664660442182604
$ws.Range("DR2").Value = -0.02493578984991634
$ws.Range("DS2").Value = -0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.06875348702283682
$ws.Range("DY2").Value = -0
$ws.Range("DZ2").Value = -0.01307614595831994
$ws.Range("EA2").Value = -0.03573607780103341
$ws.Range("EB2").Value = 0
$ws.Range("EF2").Value = -0
$ws.Range("EG2").Value = 0.04906983317781602
$ws.Range("EI2").Value = 0.08964384050486533
$ws.Range("EJ2").Value = -0.03308765702123937
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.04350483282345033
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.03706498798249654
$ws.Range("ES2").Value = 0.02057055545050679
$ws.Range("ET2").Value = 0
$ws.Range("EU2").Value = -0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.04262581418030015
$ws.Range("FA2").Value = -0.03410540786464537
$ws.Range("FB2").Value = 0.01800334742361566
$ws.Range("FD2").Value = -0
$ws.Range("FG2").Value = -0
$ws.Range("FH2").Value = -0.004411514900748431
$ws.Range("FJ2").Value = -0.01531828095525909
$ws.Range("FK2").Value = 0.008546415680021329
$ws.Range("FL2").Value = -0
$ws.Range("FN2").Value = -0
$ws.Range("FP2").Value = -0
$ws.Range("FQ2").Value = -0.01248113415308976
$ws.Range("FR2").Value = -0
$ws.Range("FS2").Value = -0.00559501352612642
$ws.Range("FT2").Value = 0.004684377672604515
$ws.Range("FV2").Value = -0
$ws.Range("FW2").Value = -0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.04212936183134562
$ws.Range("GB2").Value = 0.02691007616132334
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = -0
